$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "hv106_fctb"
$ws.Range("A11").Select()
